$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1399.78
$wsSummary.Range("B4").Value = -0.22
$wsSummary.Range("B5").Value = -2.2
$wsSummary.Range("B6").Value = 2
$wsSummary.Range("B8").Value = 2

# --- Strategy Status sheet (MarketMaking row, row 5) ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C5").Value = 99.78
$wsStatus.Range("D5").Value = 2
$wsStatus.Range("E5").Value = -0.22
$wsStatus.Range("F5").Value = -0.22

# --- New trade row data (Trade #2) ---
function Add-TradeRow($ws) {
    $ws.Range("A3").Value = 2
    $ws.Range("B3").NumberFormat = "@"
    $ws.Range("B3").Value = "2026-02-17"
    $ws.Range("C3").Value = "19:55:36"
    $ws.Range("D3").Value = "MarketMaking"
    $ws.Range("E3").Value = "UP"
    $ws.Range("F3").Value = 0.66
    $ws.Range("G3").Value = 0.59
    $ws.Range("H3").Value = "CLOSED"
    $ws.Range("I3").Value = -10.6061
    $ws.Range("J3").Value = -0.07000000000000001
    $ws.Range("K3").Value = 99.78
    $ws.Range("L3").Value = 0
    $ws.Range("M3").Value = 0
    $ws.Range("N3").Value = 0.6
    $ws.Range("O3").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P3").Value = "early_exit"
    $ws.Range("Q3").Value = 0.13
}

# --- All Trades sheet ---
$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAllTrades

# --- MarketMaking sheet ---
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $wsMarketMaking
